$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.667067050933838
$ws.Range("B1").Value = 3.658567428588867
$ws.Range("C1").Value = 2.154602289199829
$ws.Range("D1").Value = 1.505482316017151
$ws.Range("E1").Value = 1.286141872406006
